$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the header row (row 1); the data rows shift up to become rows 1-3.
$ws.Rows.Item(1).Delete()

# Hyperlink ranges aren't shifted automatically by the row delete, so
# rebuild them pointing at the correct (now one-row-higher) cells, in
# the same order and with the same targets, so the relationship ids
# still come out as rId1..rId5.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C1"), "mailto:user1@mail.com", "", "", "user1@mail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "http://site-2/", "", "", "http://site-2.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:user2@mail", "", "", "user2@mail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "http://site-3/", "", "", "http://site-3.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:user3@mail.com", "", "", "user3@mail.com")

# The Email column has always used a plain blue (non-underlined) font,
# not Excel's default underlined hyperlink look - restore that.
$ws.Range("C1").Font.Underline = -4142
$ws.Range("C2").Font.Underline = -4142
$ws.Range("C3").Font.Underline = -4142

# Row 3's Site cell (now row 2) carried a stray yellow highlight fill
# and also just picked up the default hyperlink font; clean it up by
# copying the plain look used by the rest of the Site column.
$ws.Range("B1").Copy()
$ws.Range("B2:B3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C1").Select()
